$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = [double]"0.1062659480842234"
$ws.Range("B3").Value = [double]"0.002353542544731945"
$ws.Range("C3").Value = [double]"0.0007249597814497836"
$ws.Range("D3").Value = [double]"3.401411429781534"
$ws.Range("E3").Value = [double]"0.01146848643362273"
$ws.Range("F3").Value = [double]"0.000932641161065022"
$ws.Range("G3").Value = [double]"0.003774443928398867"
$ws.Range("H3").Value = [double]"0.1086194906289554"
$ws.Range("B4").Value = [double]"0.01551897117455473"
$ws.Range("C4").Value = [double]"0.001390619828720225"
$ws.Range("D4").Value = [double]"7.916538252138063"
$ws.Range("E4").Value = [double]"0.003654032766058942"
$ws.Range("F4").Value = [double]"0.0127933958623236"
$ws.Range("G4").Value = [double]"0.01824454648678586"
$ws.Range("H4").Value = [double]"0.1217849192587782"
$ws.Range("B5").Value = [double]"0.0237162710381657"
$ws.Range("C5").Value = [double]"0.003080954633204196"
$ws.Range("D5").Value = [double]"8.924855489714673"
$ws.Range("E5").Value = [double]"0.07990755956350835"
$ws.Range("F5").Value = [double]"0.01767767470357442"
$ws.Range("G5").Value = [double]"0.02975486737275699"
$ws.Range("H5").Value = [double]"0.1299822191223892"
$ws.Range("B6").Value = [double]"0.01210807925980701"
$ws.Range("C6").Value = [double]"0.002448616276785519"
$ws.Range("D6").Value = [double]"4.003078275249077"
$ws.Range("E6").Value = [double]"0.008493214677871069"
$ws.Range("F6").Value = [double]"0.007308865147837591"
$ws.Range("G6").Value = [double]"0.01690729337177643"
$ws.Range("H6").Value = [double]"0.1183740273440305"
$ws.Range("B7").Value = [double]"0.01210602768469882"
$ws.Range("C7").Value = [double]"0.002639353696956034"
$ws.Range("D7").Value = [double]"2.540831676764028"
$ws.Range("E7").Value = [double]"0.0004651155619674327"
$ws.Range("F7").Value = [double]"0.00693297370969636"
$ws.Range("G7").Value = [double]"0.01727908165970128"
$ws.Range("H7").Value = [double]"0.1183719757689223"
$ws.Range("B8").Value = [double]"0.01340306195374568"
$ws.Range("C8").Value = [double]"0.003356760690397372"
$ws.Range("D8").Value = [double]"3.921867453231823"
$ws.Range("E8").Value = [double]"0.05956448983900994"
$ws.Range("F8").Value = [double]"0.006823911897315884"
$ws.Range("G8").Value = [double]"0.01998221201017548"
$ws.Range("H8").Value = [double]"0.1196690100379691"
$ws.Range("B9").Value = [double]"0.01755732835670036"
$ws.Range("C9").Value = [double]"0.006348964754802574"
$ws.Range("D9").Value = [double]"5.662724182121321"
$ws.Range("E9").Value = [double]"0.05818078646150379"
$ws.Range("F9").Value = [double]"0.005113548573190635"
$ws.Range("G9").Value = [double]"0.03000110814021008"
$ws.Range("H9").Value = [double]"0.1238232764409238"
$ws.Range("B10").Value = [double]"-0.1062659480842234"
$ws.Range("C10").Value = [double]"0.0005424669899068639"
$ws.Range("D10").Value = [double]"-225.1844759037409"
$ws.Range("E10").Value = [double]"0"
$ws.Range("F10").Value = [double]"-0.1073291684002668"
$ws.Range("G10").Value = [double]"-0.1052027277681801"
$ws.Range("B11").Value = [double]"-0.04715245535575163"
$ws.Range("C11").Value = [double]"0.0005847710949516904"
$ws.Range("D11").Value = [double]"-91.89482276972018"
$ws.Range("E11").Value = [double]"4.216738132124123e-147"
$ws.Range("F11").Value = [double]"-0.04829859053817288"
$ws.Range("G11").Value = [double]"-0.04600632017333038"
$ws.Range("H11").Value = [double]"0.05911349272847181"
$ws.Range("B12").Value = [double]"-0.03728263385311863"
$ws.Range("C12").Value = [double]"0.0005679261421352681"
$ws.Range("D12").Value = [double]"-75.67148455398188"
$ws.Range("E12").Value = [double]"2.583513092750158e-127"
$ws.Range("F12").Value = [double]"-0.0383957534195764"
$ws.Range("G12").Value = [double]"-0.03616951428666085"
$ws.Range("H12").Value = [double]"0.06898331423110482"
$ws.Range("B13").Value = [double]"-0.03194050580889547"
$ws.Range("C13").Value = [double]"0.0005615585538145604"
$ws.Range("D13").Value = [double]"-64.71047806309645"
$ws.Range("E13").Value = [double]"1.098501194656275e-99"
$ws.Range("F13").Value = [double]"-0.03304114507831326"
$ws.Range("G13").Value = [double]"-0.03083986653947768"
$ws.Range("H13").Value = [double]"0.07432544227532797"
$ws.Range("B14").Value = [double]"-0.02714634814911521"
$ws.Range("C14").Value = [double]"0.0005485467363789416"
$ws.Range("D14").Value = [double]"-58.58160996414945"
$ws.Range("E14").Value = [double]"1.113377849169694e-56"
$ws.Range("F14").Value = [double]"-0.02822148463471231"
$ws.Range("G14").Value = [double]"-0.02607121166351811"
$ws.Range("H14").Value = [double]"0.07911959993510824"
$ws.Range("B15").Value = [double]"-0.0238260743944529"
$ws.Range("C15").Value = [double]"0.0005433205134184146"
$ws.Range("D15").Value = [double]"-51.48557711962901"
$ws.Range("E15").Value = [double]"4.807417133721717e-33"
$ws.Range("F15").Value = [double]"-0.02489096760964443"
$ws.Range("G15").Value = [double]"-0.02276118117926137"
$ws.Range("H15").Value = [double]"0.08243987368977054"
$ws.Range("B16").Value = [double]"-0.0226544908008753"
$ws.Range("C16").Value = [double]"0.0005413481824493235"
$ws.Range("D16").Value = [double]"-49.14380947023087"
$ws.Range("E16").Value = [double]"9.184415722130954e-12"
$ws.Range("F16").Value = [double]"-0.0237155183621775"
$ws.Range("G16").Value = [double]"-0.02159346323957312"
$ws.Range("H16").Value = [double]"0.08361145728334815"
$ws.Range("B17").Value = [double]"-0.02058095543646835"
$ws.Range("C17").Value = [double]"0.0005436941324481073"
$ws.Range("D17").Value = [double]"-45.26723661290825"
$ws.Range("E17").Value = [double]"0.0007138671740878653"
$ws.Range("F17").Value = [double]"-0.02164658090522065"
$ws.Range("G17").Value = [double]"-0.01951532996771605"
$ws.Range("H17").Value = [double]"0.08568499264775509"
$ws.Range("B18").Value = [double]"-0.01872039579824159"
$ws.Range("C18").Value = [double]"0.0005488260628087845"
$ws.Range("D18").Value = [double]"-39.43301199692144"
$ws.Range("E18").Value = [double]"1.43470515300722e-17"
$ws.Range("F18").Value = [double]"-0.01979607973422014"
$ws.Range("G18").Value = [double]"-0.01764471186226304"
$ws.Range("H18").Value = [double]"0.08754555228598186"
$ws.Range("B19").Value = [double]"-0.01513411775007059"
$ws.Range("C19").Value = [double]"0.0005471785277015311"
$ws.Range("D19").Value = [double]"-32.13117341373442"
$ws.Range("E19").Value = [double]"4.734724551759598e-07"
$ws.Range("F19").Value = [double]"-0.01620657258284563"
$ws.Range("G19").Value = [double]"-0.01406166291729555"
$ws.Range("H19").Value = [double]"0.09113183033415286"
$ws.Range("B20").Value = [double]"-0.01193545640062883"
$ws.Range("C20").Value = [double]"0.0005580973809477267"
$ws.Range("D20").Value = [double]"-23.15374294766456"
$ws.Range("E20").Value = [double]"0.000137991406657008"
$ws.Range("F20").Value = [double]"-0.01302931189854052"
$ws.Range("G20").Value = [double]"-0.01084160090271713"
$ws.Range("H20").Value = [double]"0.09433049168359461"
$ws.Range("B21").Value = [double]"-0.009207007581023979"
$ws.Range("C21").Value = [double]"0.0005684379887194317"
$ws.Range("D21").Value = [double]"-18.13655578945873"
$ws.Range("E21").Value = [double]"0.08992082778524013"
$ws.Range("F21").Value = [double]"-0.01032113043782094"
$ws.Range("G21").Value = [double]"-0.008092884724227012"
$ws.Range("H21").Value = [double]"0.09705894050319946"
$ws.Range("B22").Value = [double]"-0.006422585835077638"
$ws.Range("C22").Value = [double]"0.0005674470431213879"
$ws.Range("D22").Value = [double]"-11.41648966504769"
$ws.Range("E22").Value = [double]"0.05022836113010615"
$ws.Range("F22").Value = [double]"-0.007534766535506167"
$ws.Range("G22").Value = [double]"-0.005310405134649107"
$ws.Range("H22").Value = [double]"0.09984336224914581"
$ws.Range("B23").Value = [double]"-0.004742375664289536"
$ws.Range("C23").Value = [double]"0.0005777206480378481"
$ws.Range("D23").Value = [double]"-8.510670225173843"
$ws.Range("E23").Value = [double]"0.1174109724983757"
$ws.Range("F23").Value = [double]"-0.005874692339873002"
$ws.Range("G23").Value = [double]"-0.003610058988706069"
$ws.Range("H23").Value = [double]"0.1015235724199339"
$ws.Range("B24").Value = [double]"-0.003446680604358203"
$ws.Range("C24").Value = [double]"0.0005676707874630037"
$ws.Range("D24").Value = [double]"-6.219936886009107"
$ws.Range("E24").Value = [double]"0.1089460705801795"
$ws.Range("F24").Value = [double]"-0.004559299798456018"
$ws.Range("G24").Value = [double]"-0.002334061410260388"
$ws.Range("H24").Value = [double]"0.1028192674798653"
$ws.Range("B25").Value = [double]"-0.002719240918236039"
$ws.Range("C25").Value = [double]"0.0005524611953109567"
$ws.Range("D25").Value = [double]"-5.39113066075245"
$ws.Range("E25").Value = [double]"0.2293298616602996"
$ws.Range("F25").Value = [double]"-0.003802049681633672"
$ws.Range("G25").Value = [double]"-0.001636432154838407"
$ws.Range("H25").Value = [double]"0.1035467071659874"
$ws.Range("B26").Value = [double]"0.01915803702108339"
$ws.Range("C26").Value = [double]"0.002494260875470018"
$ws.Range("D26").Value = [double]"17.40871393610974"
$ws.Range("E26").Value = [double]"0.008179815621105877"
$ws.Range("F26").Value = [double]"0.01426935967355216"
$ws.Range("G26").Value = [double]"0.02404671436861464"
$ws.Range("H26").Value = [double]"0.1254239851053068"
